# Update the "Overview" sheet of the Feghadir (folad) yearly income-statement
# workbook: roll the five reporting periods forward by one year
# (1396-1400 -> 1397-1401), refresh the "publish date" row, and replace the
# financial figures with the newly published numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- Row 8: financial-period column headers (shift one year forward) ---
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# --- Row 9: publish-date row (shift one year forward) ---
$ws.Range("D9").Value = "1400-09-29 (2)"
$ws.Range("E9").Value = "1400-10-25 (4)"
$ws.Range("F9").Value = "1401-03-07 (8)"
$ws.Range("G9").Value = "1402-03-11 (9)"
$ws.Range("H9").Value = "1402-03-11 (2)"

# --- Row 11: فروش (Sales) ---
$ws.Range("D11").Value = 107727
$ws.Range("E11").Value = 131732
$ws.Range("F11").Value = 151620
$ws.Range("G11").Value = 290249
$ws.Range("H11").Value = 190316

# --- Row 12: بهای تمام شده کالای فروش رفته (Cost of goods sold) ---
$ws.Range("D12").Value = -72923
$ws.Range("E12").Value = -103031
$ws.Range("F12").Value = -92294
$ws.Range("G12").Value = -203682
$ws.Range("H12").Value = -129297

# --- Row 13: سود (زیان) ناخالص (Gross profit) ---
$ws.Range("D13").Value = 34804
$ws.Range("E13").Value = 28701
$ws.Range("F13").Value = 59326
$ws.Range("G13").Value = 86567
$ws.Range("H13").Value = 61018

# --- Row 14: هزینه های عمومی, اداری و تشکیلاتی (G&A expenses) ---
$ws.Range("D14").Value = -784
$ws.Range("E14").Value = -2735
$ws.Range("F14").Value = -2117
$ws.Range("G14").Value = -865
$ws.Range("H14").Value = -1543

# --- Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی (Other operating income/expense, net) ---
$ws.Range("D16").Value = 1198
$ws.Range("E16").Value = 501
$ws.Range("F16").Value = 806
$ws.Range("G16").Value = -910
$ws.Range("H16").Value = 3327

# --- Row 17: سود (زیان) عملیاتی (Operating profit) ---
$ws.Range("D17").Value = 35218
$ws.Range("E17").Value = 26467
$ws.Range("F17").Value = 58014
$ws.Range("G17").Value = 84793
$ws.Range("H17").Value = 62802

# --- Row 18: هزینه های مالی (Financial expenses) ---
$ws.Range("D18").Value = -12
$ws.Range("E18").Value = -10
$ws.Range("F18").Value = -28
$ws.Range("G18").Value = -242
$ws.Range("H18").Value = -823

# --- Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی (Other non-operating income/expense, net) ---
$ws.Range("D19").Value = 5857
$ws.Range("E19").Value = 3833
$ws.Range("F19").Value = 6653
$ws.Range("G19").Value = 12196
$ws.Range("H19").Value = 9731

# --- Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات (Pre-tax profit from continuing operations) ---
$ws.Range("D20").Value = 41063
$ws.Range("E20").Value = 30290
$ws.Range("F20").Value = 64639
$ws.Range("G20").Value = 96746
$ws.Range("H20").Value = 71710

# --- Row 21: مالیات (Tax) ---
$ws.Range("D21").Value = -83
$ws.Range("E21").Value = -192
$ws.Range("F21").Value = -46
$ws.Range("G21").Value = -115
$ws.Range("H21").Value = -218

# --- Row 22: سود (زیان) خالص عملیات در حال تداوم (Net profit from continuing operations) ---
$ws.Range("D22").Value = 40980
$ws.Range("E22").Value = 30098
$ws.Range("F22").Value = 64593
$ws.Range("G22").Value = 96631
$ws.Range("H22").Value = 71491

# --- Row 24: سود (زیان) خالص (Net profit) ---
$ws.Range("D24").Value = 40980
$ws.Range("E24").Value = 30098
$ws.Range("F24").Value = 64593
$ws.Range("G24").Value = 96631
$ws.Range("H24").Value = 71491

# --- Row 26: سرمایه (Capital) ---
$ws.Range("D26").Value = 44485
$ws.Range("E26").Value = 35079
$ws.Range("F26").Value = 53073
$ws.Range("G26").Value = 45480
$ws.Range("H26").Value = 34005
